$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1. Insert blank column before column E (shifts coeff_s..eta_d and beyond to the right)
$ws.Columns.Item(5).Insert()

# 2. Set new header name first (so the engine treats E1 write as just a value set, not yet resize)
$ws.Range("E1").Value = "Ref int ratio"

# 3. Resize table to include the new column
$newRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(11,11))
$lo.Resize($newRange)

# 4. Restore names for the shifted columns (F..K) back to their original names
$ws.Range("F1").Value = "coeff_s"
$ws.Range("G1").Value = "coeff_p"
$ws.Range("H1").Value = "coeff_d"
$ws.Range("I1").Value = "eta_s"
$ws.Range("J1").Value = "eta_p"
$ws.Range("K1").Value = "eta_d"

Write-Host "Range: $($lo.Range.Address())"
for ($i = 1; $i -le $lo.ListColumns.Count; $i++) {
    $col = $lo.ListColumns.Item($i)
    Write-Host "$i : $($col.Name)"
}

# 5. Set the calculated column formula for the new "Ref int ratio" column
$ws.Range("E2:E11").Formula = "=Table1[[#This Row],[Int ratio]]"
Write-Host "----"
Write-Host $ws.Range("E2").Formula
Write-Host $ws.Range("E2").Value
